# Auto-generated: applies scheduled market-price/profit refresh values
# to the Mateus_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2500969.2
$ws.Range("J17").Value = 2500969.2
$ws.Range("L17").Value = 7502907.600000001
$ws.Range("N17").Value = -7503243.600000001
$ws.Range("H48").Value = 2594.3333
$ws.Range("J48").Value = 5833
$ws.Range("L48").Value = 17499
$ws.Range("N48").Value = -18083
$ws.Range("H56").Value = 2594.3333
$ws.Range("J56").Value = 5833
$ws.Range("L56").Value = 17499
$ws.Range("N56").Value = -18567
$ws.Range("H64").Value = 10222.125
$ws.Range("I64").Value = 3311
$ws.Range("K64").Value = 3311
$ws.Range("M64").Value = -3063
$ws.Range("H67").Value = 10222.125
$ws.Range("I67").Value = 3311
$ws.Range("K67").Value = 3311
$ws.Range("M67").Value = -2453
$ws.Range("H70").Value = 4003.8076
$ws.Range("I70").Value = 1700.25
$ws.Range("K70").Value = 5100.75
$ws.Range("M70").Value = -4830.75
$ws.Range("H73").Value = 4003.8076
$ws.Range("I73").Value = 1700.25
$ws.Range("K73").Value = 5100.75
$ws.Range("M73").Value = -4164.75
$ws.Range("H100").Value = 1916
$ws.Range("J100").Value = 3166.6667
$ws.Range("L100").Value = 3166.6667
$ws.Range("N100").Value = -4248.6667

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 18000
$ws.Range("I37").Value = 18000
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 18000
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -17727
$ws.Range("N37").ClearContents()  # was -21546
$ws.Range("H45").Value = 5839.2666
$ws.Range("I45").Value = 2600
$ws.Range("J45").Value = 7998.778
$ws.Range("K45").Value = 2600
$ws.Range("L45").Value = 7998.778
$ws.Range("M45").Value = -2223
$ws.Range("N45").Value = -8752.778
$ws.Range("H63").Value = 6016.6
$ws.Range("I63").Value = 3833
$ws.Range("K63").Value = 3833
$ws.Range("M63").Value = -3147
$ws.Range("H66").Value = 6016.6
$ws.Range("I66").Value = 3833
$ws.Range("K66").Value = 19165
$ws.Range("M66").Value = -15733
$ws.Range("H74").Value = 3869.5862
$ws.Range("I74").Value = 3237.7896
$ws.Range("K74").Value = 3237.7896
$ws.Range("M74").Value = -2363.7896
$ws.Range("H77").Value = 3869.5862
$ws.Range("I77").Value = 3237.7896
$ws.Range("K77").Value = 16188.948
$ws.Range("M77").Value = -11820.948
$ws.Range("H122").Value = 1602.3684
$ws.Range("I122").Value = 1292.7142
$ws.Range("K122").Value = 3878.1426
$ws.Range("M122").Value = -1428.1426
$ws.Range("H132").Value = 5181.528
$ws.Range("I132").Value = 4271.6772
$ws.Range("J132").Value = 10822.6
$ws.Range("K132").Value = 12815.0316
$ws.Range("L132").Value = 32467.8
$ws.Range("M132").Value = -10285.0316
$ws.Range("N132").Value = -37527.8

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()  # was -2940
$ws.Range("H31").Value = 5145.1304
$ws.Range("J31").Value = 5410.846
$ws.Range("L31").Value = 5410.846
$ws.Range("N31").Value = -6000.846
$ws.Range("H34").Value = 5145.1304
$ws.Range("J34").Value = 5410.846
$ws.Range("L34").Value = 5410.846
$ws.Range("N34").Value = -5814.846
$ws.Range("H58").Value = 5447.4814
$ws.Range("I58").Value = 3004.611
$ws.Range("K58").Value = 3004.611
$ws.Range("M58").Value = -2801.611
$ws.Range("H111").Value = 79973
$ws.Range("J111").Value = 79973
$ws.Range("L111").Value = 79973
$ws.Range("N111").Value = -88153
$ws.Range("H132").Value = 2373.158
$ws.Range("I132").Value = 2012.2858
$ws.Range("K132").Value = 6036.857400000001
$ws.Range("M132").Value = -3506.857400000001
$ws.Range("H134").Value = 5548.64
$ws.Range("I134").Value = 3829.45
$ws.Range("K134").Value = 11488.35
$ws.Range("M134").Value = -8953.349999999999
$ws.Range("H136").Value = 5447.4814
$ws.Range("I136").Value = 3004.611
$ws.Range("K136").Value = 9013.832999999999
$ws.Range("M136").Value = -6463.832999999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 3224.25
$ws.Range("J12").Value = 3224.25
$ws.Range("L12").Value = 9672.75
$ws.Range("N12").Value = -10018.75
$ws.Range("H75").Value = 3666.6667
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 3666.6667
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 11000.0001
$ws.Range("M75").ClearContents()  # was 548
$ws.Range("N75").Value = -12996.0001
$ws.Range("H78").Value = 3666.6667
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 3666.6667
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 33000.0003
$ws.Range("M78").ClearContents()  # was 3642
$ws.Range("N78").Value = -42984.0003
$ws.Range("H98").Value = 317.66666
$ws.Range("I98").Value = 326.5
$ws.Range("K98").Value = 979.5
$ws.Range("M98").Value = 518.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 199
$ws.Range("I5").Value = 199
$ws.Range("K5").Value = 199
$ws.Range("M5").Value = -87
$ws.Range("H9").Value = 732.2
$ws.Range("I9").Value = 911.5
$ws.Range("J9").Value = 15
$ws.Range("K9").Value = 911.5
$ws.Range("L9").Value = 15
$ws.Range("M9").Value = -741.5
$ws.Range("N9").Value = -355
$ws.Range("H134").Value = 59994.2
$ws.Range("J134").Value = 59994.2
$ws.Range("L134").Value = 179982.6
$ws.Range("N134").Value = -185052.6

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2319.6
$ws.Range("I7").Value = 2319.6
$ws.Range("K7").Value = 2319.6
$ws.Range("M7").Value = -2207.6
$ws.Range("H9").Value = 3375.4
$ws.Range("I9").Value = 3375.4
$ws.Range("K9").Value = 3375.4
$ws.Range("M9").Value = -3151.4
$ws.Range("H29").Value = 5000
$ws.Range("I29").Value = 5000
$ws.Range("K29").Value = 5000
$ws.Range("M29").Value = -4705
$ws.Range("H30").Value = 9899.5
$ws.Range("I30").Value = 9899.5
$ws.Range("K30").Value = 9899.5
$ws.Range("M30").Value = -9791.5
$ws.Range("H68").Value = 2350
$ws.Range("I68").Value = 2325
$ws.Range("J68").Value = 2400
$ws.Range("K68").Value = 2325
$ws.Range("L68").Value = 2400
$ws.Range("M68").Value = -1576
$ws.Range("N68").Value = -3898
$ws.Range("H71").Value = 2350
$ws.Range("I71").Value = 2325
$ws.Range("J71").Value = 2400
$ws.Range("K71").Value = 11625
$ws.Range("L71").Value = 12000
$ws.Range("M71").Value = -7881
$ws.Range("N71").Value = -19488
$ws.Range("H126").Value = 2319.6
$ws.Range("I126").Value = 2319.6
$ws.Range("K126").Value = 6958.799999999999
$ws.Range("M126").Value = -4488.799999999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4360.5
$ws.Range("I62").Value = 4430.625
$ws.Range("J62").Value = 4267
$ws.Range("K62").Value = 4430.625
$ws.Range("L62").Value = 4267
$ws.Range("M62").Value = -3806.625
$ws.Range("N62").Value = -5515
$ws.Range("H65").Value = 4360.5
$ws.Range("I65").Value = 4430.625
$ws.Range("J65").Value = 4267
$ws.Range("K65").Value = 22153.125
$ws.Range("L65").Value = 21335
$ws.Range("M65").Value = -19033.125
$ws.Range("N65").Value = -27575
$ws.Range("H132").Value = 2942.7441
$ws.Range("I132").Value = 2828
$ws.Range("K132").Value = 8484
$ws.Range("M132").Value = -6954

